$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add Trust_Name / RR8 columns and relabel existing headers
$ws.Range("C1").Value = "Trust_Name"
$ws.Range("C3").Value = "RR8"
$ws.Range("B1").Value = "Initial Label"
$ws.Range("D1").Value = "Expected Graph Label"
$ws.Range("B3").Value = "Region R1 and Others"
$ws.Range("C4").Value = "RR8"
$ws.Range("B4").Value = "Region R1 and Others"
$ws.Range("D4").Value = "Region R1 and Others"
$ws.Range("C5").Value = "RR8"
$ws.Range("B5").Value = "Region R1 and Others"
$ws.Range("D5").Value = "Region R1 and Others"

# Re-assert UserName/Password/admin/others block so it keeps its original shared-string slots
$ws.Range("E2").Value = "admin"
$ws.Range("F2").Value = "admin"

# New thin left/right border around the two "RR8" trust-name helper cells
$ws.Range("D1").Borders.Item(7).LineStyle = 1
$ws.Range("D1").Borders.Item(7).ColorIndex = -4105
$ws.Range("D1").Borders.Item(10).LineStyle = 1
$ws.Range("D1").Borders.Item(10).ColorIndex = -4105

$ws.Range("C5").Borders.Item(7).LineStyle = 1
$ws.Range("C5").Borders.Item(7).ColorIndex = -4105
$ws.Range("C5").Borders.Item(10).LineStyle = 1
$ws.Range("C5").Borders.Item(10).ColorIndex = -4105

# New column widths for the inserted C/D columns
$ws.Columns("C").ColumnWidth = 10.8333333333333
$ws.Columns("D").ColumnWidth = 19.5

# Move the active selection
$ws.Range("E3").Select()
